$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Power FET changed from SIR5802DP-T1-RE3 to IRFB7730PBF (TO-220 package)
$ws.Range("A3").Value = "IRFB7730PBF-ND"
$ws.Range("B3").Value = "IRFB7730PBF"
$ws.Range("E3").Value = 3.25

# Update the active selection to A24, matching the author's final cursor position
$ws.Range("A24").Select()
